# table_device_package.xlsx - "add some info as commentary add Xiaomi Support"
#
# Adds two new IntelliJ/Darcula-style "code snippet" rich-text cells (C8, C9)
# documenting the MIUI_ACTION / MIUI_EXTRA constants used for Xiaomi's
# "App permission editor" intent, matching the existing code-snippet styling
# already used elsewhere in the sheet (e.g. E7/C6), and nudges the selection
# to where the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# C8: `private static final String MIUI_ACTION = "miui.intent.action.APP_PERM_EDITOR";`
#   - base cell font/color = Courier New 10pt FFCC7832 (keyword/punctuation
#     color) -> covers the leading "private static final " run (no override)
#     and the trailing ";" run.
#   - "String " and "= " -> FFA9B7C6 (type/operator color)
#   - "MIUI_ACTION " -> italic FF9876AA (identifier color)
#   - the quoted literal -> FF6A8759 (string-literal color)
# ---------------------------------------------------------------------------
$ws.Range("C8").Value = "private static final String MIUI_ACTION = ""miui.intent.action.APP_PERM_EDITOR"";"
$c8 = $ws.Range("C8")
$c8.Font.Name = "Courier New"
$c8.Font.Size = 10
$c8.Font.Color = 3307724
$c8.VerticalAlignment = -4108
$c8.Characters(22, 7).Font.Name = "Courier New"
$c8.Characters(22, 7).Font.Size = 10
$c8.Characters(22, 7).Font.Color = 13023145
$c8.Characters(29, 12).Font.Name = "Courier New"
$c8.Characters(29, 12).Font.Size = 10
$c8.Characters(29, 12).Font.Color = 11171480
$c8.Characters(29, 12).Font.Italic = $true
$c8.Characters(41, 2).Font.Name = "Courier New"
$c8.Characters(41, 2).Font.Size = 10
$c8.Characters(41, 2).Font.Color = 13023145
$c8.Characters(43, 36).Font.Name = "Courier New"
$c8.Characters(43, 36).Font.Size = 10
$c8.Characters(43, 36).Font.Color = 5867370
$c8.Characters(79, 1).Font.Name = "Courier New"
$c8.Characters(79, 1).Font.Size = 10
$c8.Characters(79, 1).Font.Color = 3307724

# ---------------------------------------------------------------------------
# C9: `private static final String MIUI_EXTRA = "extra_pkgname";`
#   same color scheme as C8, plus the dark "code block" background fill
#   (FF2B2B2B) already used by the other highlighted code cells.
# ---------------------------------------------------------------------------
$ws.Range("C9").Value = "private static final String MIUI_EXTRA = ""extra_pkgname"";"
$c9 = $ws.Range("C9")
$c9.Font.Name = "Courier New"
$c9.Font.Size = 10
$c9.Font.Color = 3307724
$c9.Interior.Color = 2829099
$c9.VerticalAlignment = -4108
$c9.Characters(22, 7).Font.Name = "Courier New"
$c9.Characters(22, 7).Font.Size = 10
$c9.Characters(22, 7).Font.Color = 13023145
$c9.Characters(29, 11).Font.Name = "Courier New"
$c9.Characters(29, 11).Font.Size = 10
$c9.Characters(29, 11).Font.Color = 11171480
$c9.Characters(29, 11).Font.Italic = $true
$c9.Characters(40, 2).Font.Name = "Courier New"
$c9.Characters(40, 2).Font.Size = 10
$c9.Characters(40, 2).Font.Color = 13023145
$c9.Characters(42, 15).Font.Name = "Courier New"
$c9.Characters(42, 15).Font.Size = 10
$c9.Characters(42, 15).Font.Color = 5867370
$c9.Characters(57, 1).Font.Name = "Courier New"
$c9.Characters(57, 1).Font.Size = 10
$c9.Characters(57, 1).Font.Color = 3307724

# ---------------------------------------------------------------------------
# Row heights grew slightly (author's Excel re-measured wrapped-text rows);
# reproduce the new explicit heights for the affected rows.
# ---------------------------------------------------------------------------
$ws.Rows(3).RowHeight = 105
$ws.Rows(5).RowHeight = 45
$ws.Rows(6).RowHeight = 60
$ws.Rows(8).RowHeight = 30

# ---------------------------------------------------------------------------
# Selection moved to E7 by the time the author saved.
# ---------------------------------------------------------------------------
$ws.Range("E7").Select() | Out-Null
